$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the Sep 25 2023 GitHub Actions cryptos-list refresh.
# Percent-change cells keep their padded "  +x.xx%  " / "  -x.xx%  " text layout.
# Price cells that are plain decimals (e.g. "2.17") would otherwise be auto-
# detected as numbers by Excel, so those are briefly marked as Text (@) before
# the write and restored to the Normal style afterwards to avoid changing the
# cell formatting/style - only the underlying text value changes.

$ws.Range('D2').Value = '26.175.01'
$ws.Range('E2').Value = '  -2.11%  '
$ws.Range('D3').Value = '1.576.41'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  -0.49%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '208.87'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('E6').Value = '  -2.89%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  -1.57%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.245'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -1.16%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '19.58'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = '1.798.28'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.06'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.577.37'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('E15').Value = '  -1.87%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '64.43'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '26.159.46'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('E18').Value = '  -1.91%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.26'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.67%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '208.66'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('E21').Value = '  -0.40%  '
$ws.Range('E22').Value = '  -0.94%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.45%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '8.84'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.39%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '143.77'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('E26').Value = '  -0.47%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '6.99'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('E28').Value = '  -1.40%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '15.21'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  -1.39%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.21'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.98%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.01'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').Value = '1.279.40'
$ws.Range('E34').Value = '  -0.59%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.613'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +4.30%  '
$ws.Range('E36').Value = '  -1.46%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.48'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -9.79%  '
$ws.Range('E39').Value = '  -2.35%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.811'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('E41').Value = '  -0.41%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.59'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +2.51%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.764'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('E44').Value = '  -2.93%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '62.39'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('D46').Value = '1.711.00'
$ws.Range('E46').Value = '  -1.66%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '88.78'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E48').Value = '  -2.31%  '
$ws.Range('E49').Value = '  -4.05%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.101'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('E51').Value = '  -1.59%  '
